$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRevisions = @{
  2 = 0.000000000000009992007221626409
  3 = -0.000000000000009547918011776346
  4 = -0.000000000000004440892098500626
  5 = 0.000000000000008160139230994901
  6 = 0.000000000000002997602166487923
  7 = 0.00000000000001021405182655144
  8 = -0.00000000000002442490654175344
  9 = 0.00000000000002436939539052219
  10 = -0.00000000000001865174681370263
  11 = -0.000000000000009103828801926284
  12 = 0.00000000000002264854970235319
  13 = 1.593775310996169
  14 = -0.4433754434956397
  15 = -0.5659094841436607
  16 = -0.02440516240001389
  17 = 0.5808185939316554
  18 = 0.08546502771525721
  19 = 0.1413663207244666
  20 = -0.9925327920106826
  21 = 0.6166013166907478
  22 = 0.1324346314480873
  23 = -0.1095386862888253
  24 = -0.2383424898041266
  25 = 0.6794258297128992
  26 = 0.03106364948285656
  27 = -0.4886483869543554
  28 = 0.605758692296735
  29 = 0.2658730860130905
  30 = -0.215129482578987
  31 = -0.4485771038079694
  32 = 0.3419104562980941
  33 = 0.1184663661808556
  34 = 0.4226346167906511
  35 = -0.8394267946364184
  36 = 0.9539990704374277
  37 = -0.1016472839247859
  38 = 0.3894274519612926
  39 = -0.4632625443669942
  40 = -0.005792791899517091
  41 = 0.5624759907254965
  42 = 0.3796720383020056
  43 = -0.3492652477906187
  44 = -0.715889171321582
  45 = 0.1691567084338567
  46 = -0.6525458734825007
  47 = 0.2152667692021368
  48 = -0.885744281990282
  49 = 0.3505440551774829
  50 = 0.4918016298087824
  51 = 0.6329113189306338
  52 = -1.755361674695274
  53 = 0.2210295811551402
  54 = -0.1411005845595663
  55 = 0.05762549224104685
  56 = 0.639181144805368
  57 = -0.1144610931020827
  58 = -0.4184709022183099
  59 = 0.2429549719812211
  60 = 1.128237710175259
  61 = -0.2082144774214872
  62 = -0.6910508056900682
  63 = 0.0947112025147417
  64 = 0.5952133710329237
  65 = -0.1142735084420681
  66 = 0.4225165845712571
  67 = -0.1967006994430545
  68 = 0.8880062631184649
  69 = -0.5897473044240461
  70 = -0.1661013782684851
  71 = -0.3687592045110376
  72 = -0.7537458141733666
  73 = -0.8114147307874704
  74 = -0.1212771608465997
  75 = -0.5935793456314062
  76 = 0.3762390903770712
  77 = -0.4008196643353814
  78 = 0.1721949556597677
  79 = -0.4748495211968715
  80 = 0.9462663862015542
  81 = -0.1356194181296091
  82 = 0.02323282939418525
  83 = -0.5550776869004674
  84 = 0.8865631931960101
  85 = 0.07104253540100011
  86 = 0.1075028472288783
  87 = -0.758469243884506
  88 = -0.5661311224158858
  89 = 0.01540602510147604
  95 = -0.1711104671482269
  96 = 0.2949881563205918
  97 = 0.1633555277145124
  98 = 0.03022793029104015
  99 = -0.2928552586432669
  100 = 0.6214648934017475
  101 = -0.4443674518509457
  102 = -0.1898431833400273
  103 = 0.2351372936899813
  104 = -0.2314726160725444
  105 = -0.1398984511979687
  106 = -0.4898424607665044
  107 = -0.0619887154334321
  108 = 0.4897593076320608
  109 = -0.08025008101715425
  110 = -0.2209309121121734
  111 = -0.2447597831269769
  112 = -0.213315487371266
  113 = 0.188451994675966
  114 = -0.1625432761390304
  115 = -0.03627262647648394
  116 = -0.6870776560732526
  117 = -0.09599133664528381
  118 = 0.02834167208306626
  119 = -0.4067771049315971
  120 = 0.8473145294465119
  121 = -0.3343650428035478
  122 = 0.4558229525073971
  123 = -0.4284561131352164
  124 = -0.2317013822699299
  125 = -0.1321993960636193
  126 = -0.2990191236833212
  127 = 0.3619968567220385
  128 = -0.1919337072516534
  129 = -0.8168372431843682
  130 = -0.1962295261482101
  131 = -0.6265044734712123
  132 = -1.191162284383693
  133 = -0.7174192902096188
  134 = 1.607786042869044
  135 = -0.8848240492139481
  136 = -0.4425935756942813
  137 = -0.01202960564142147
  138 = 0.1065604846180273
  139 = -0.08591989013840001
  140 = 0.1385721760938039
  141 = 0.09472847613588289
  142 = -0.1296176279974082
  143 = -0.01074155887864159
  144 = 0.3186980753357052
  145 = 0.1914876003089772
  146 = 0.08603368373087023
  147 = -0.3828165493744078
  148 = 0.1054308279183608
}

foreach ($row in $newRevisions.Keys) {
  $ws.Cells.Item($row, 2).Value2 = $newRevisions[$row]
}

# Add new row 149 (2025-04-01), copying the date-format style from row 148
$srcFmt = $ws.Range("A148:B148")
$dstFmt = $ws.Range("A149:B149")
$srcFmt.Copy()
$dstFmt.PasteSpecial(-4122)
$ws.Range("A149").Value2 = 45748
$ws.Range("B149").Value2 = 0

Write-Host "done"
